$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# A shared date for the three new log entries (15 Aug 2025, serial 45884)
$logDate = (Get-Date -Year 2025 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0).Date

# ---------------------------------------------------------------------------
# Populate the new text in the exact order the cells were originally typed
# in, so that freshly-introduced shared strings line up the same way.
# ---------------------------------------------------------------------------
$ws.Range("C48").Value = "BFS, keep set for seen levels, if not in set append to result array"
$ws.Range("N48").Value = "Good BFS practice"
$ws.Range("C49").Value = "Map that tracks members and children of every level, and then process that map"
$ws.Range("E48").Value = "Linear (Queue is linear worst case scenario)"
$ws.Range("C50").Value = "Almost correct, but did some repetitive work"
$ws.Range("I50").Value = "I did keep track of current level in an array, and appended as level changed, but did some extra work"
$ws.Range("K50").Value = "Need to review"
$ws.Range("N50").Value = "Almost correct, perhaps not in the best mood to do the problem"

# ---------------------------------------------------------------------------
# Row 48 - "Left view"
# ---------------------------------------------------------------------------
$ws.Range("D48").Value = "Linear (visit every node once)"
$ws.Range("F48").Value = $logDate
$ws.Range("G48").Value = "10 minutes"
$ws.Range("H48").Value = "10 minutes"
$ws.Range("I48").Value = "All good"
$ws.Range("O48").Value = "No"
$ws.Range("P48").Value = 4
$ws.Range("Q48").Value = 4
$ws.Range("R48").Value = 4
$ws.Range("S48").Value = 4
$ws.Range("I48:M48").Merge()
$ws.Rows.Item(48).RowHeight = 100

# ---------------------------------------------------------------------------
# Row 49 - "Most prolific level"
# ---------------------------------------------------------------------------
$ws.Range("D49").Value = "Linear"
$ws.Range("E49").Value = "Linear"
$ws.Range("F49").Value = $logDate
$ws.Range("G49").Value = "15 minutes"
$ws.Range("H49").Value = "20 minutes"
$ws.Range("I49").Value = "All good"
$ws.Range("N49").Value = "Good BFS practice"
$ws.Range("O49").Value = "No"
$ws.Range("P49").Value = 4
$ws.Range("Q49").Value = 3
$ws.Range("R49").Value = 4
$ws.Range("S49").Value = 4
$ws.Range("I49:M49").Merge()
$ws.Rows.Item(49).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 50 - "Zig-zag order" (columns I:M stay as individual cells here)
# ---------------------------------------------------------------------------
$ws.Range("D50").Value = "Linear"
$ws.Range("E50").Value = "Linear"
$ws.Range("F50").Value = $logDate
$ws.Range("G50").Value = "15 minutes"
$ws.Range("H50").Value = "20 minutes"
$ws.Range("J50").Value = "No"
$ws.Range("L50").Value = "No"
$ws.Range("M50").Value = "No"
$ws.Range("O50").Value = "No"
$ws.Range("P50").Value = 2
$ws.Range("Q50").Value = 2
$ws.Range("R50").Value = 2
$ws.Range("S50").Value = 2
$ws.Rows.Item(50).RowHeight = 60

# ---------------------------------------------------------------------------
# Update the view: frozen pane now starts at row 47, selection moved to T50
# ---------------------------------------------------------------------------
$ws.Range("T50").Select()
$excel.ActiveWindow.ScrollRow = 47
